# Correction d'un bug qui créait des doublons de la valeur "Toutes" pour le
# choix d'années (liste IDs reçus) + statut Régulier qui ne doit pas passer à
# P&R si la somme du don régulier ne dépasse pas le minimum.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ligne 2 : GAY Noël - retire le doublon d'ID reçu et repasse "Régulier" à FAUX
$ws.Range("D2").Value = "240124NG1"
$ws.Range("G2").Value = $false

# Ligne 3 : AUBERT Nicolas - retire le doublon d'ID reçu et repasse "Régulier" à FAUX
$ws.Range("D3").Value = "240122NA1"
$ws.Range("G3").Value = $false

# Ligne 15 : PASTIER Olivia - le reçu régulier n'existe plus (don < min) donc
# on vide l'ID reçu et on remet "Régulier" à vide (plus un booléen)
$ws.Range("D15").Value = ""
$ws.Range("G15").Value = ""

# Ligne 25 : MORGAN Moral
$ws.Range("D25").Value = ""
$ws.Range("G25").Value = ""

# Ligne 37 : MARTINEAU Bérangère
$ws.Range("D37").Value = ""
$ws.Range("G37").Value = ""

# Ligne 38 : JACQUET Jacques
$ws.Range("D38").Value = ""
$ws.Range("G38").Value = ""

# Ligne 49 : WEBER Therese
$ws.Range("D49").Value = ""
$ws.Range("G49").Value = ""

# Ligne 53 : BOURDON Julien
$ws.Range("D53").Value = ""
$ws.Range("G53").Value = ""

# Ligne 60 : total sans reçus remis à 0
$ws.Range("B60").Value = 0
